$p = $ppt.ActivePresentation

# --- Slide 1: remove the leftover "Find these slides in raven/doc/misc" textbox ---
$s1 = $p.Slides.Item(1)
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -like "Find these slides in raven/doc/misc*") {
        $sh.Delete()
    }
}

# --- Slide 16: split the closing line, add a note about 'time' in dataobjects ---
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(2)
$tr = $sh16.TextFrame.TextRange

# Paragraph 7 currently reads:
#   "From inside RAVEN can access on the data object directly:	"
# Split off the trailing ":" into its own run and drop the stray tab.
$lastPara = $tr.Paragraphs(7, 1)
$run1 = $lastPara.Runs(1, 1)
$run1.Text = "From inside RAVEN can access on the data object directly"
$run1.InsertAfter(":") | Out-Null

# Append: blank paragraph, the new note (two runs), and a final tab-only paragraph.
$tr.InsertAfter("`r`rtime should not be listed as input or output variables in the new dataobjects`r`t") | Out-Null

# Fix up the wording / run split for the new "time" paragraph (paragraph 9).
$p9 = $tr.Paragraphs(9, 1)
$r9 = $p9.Runs(1, 1)
$r9.Text = "‘time’ should not be listed as input or output variables in the new "
$r9.InsertAfter("dataobjects") | Out-Null

# Final paragraph (the lone tab) should not carry a bullet.
$p10 = $tr.Paragraphs(10, 1)
$p10.ParagraphFormat.Bullet.Type = 0
